# simulator_scenarios.xlsx — add network-setup columns (J..R) in between the
# existing "young_percentage" column and the trailing infection-probability
# columns, which get pushed out to the right (old J,K,L -> new S,T,U).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert 9 blank columns before the old J (index 10), shifting
#    J,K,L (infection_prob / reinfection_prob / vaccinated_infection_prob)
#    to S,T,U. This is a real column insert (shift-right), so the existing
#    formatting/widths of J,K,L travel with the data to S,T,U untouched.
$ws.Range($ws.Columns.Item(10), $ws.Columns.Item(18)).Insert(-4161)

# 2) New header row (row 1) labels for the freshly inserted columns.
$ws.Range("J1").Value = "network_type"
$ws.Range("K1").Value = "network_param_k"
$ws.Range("L1").Value = "network_param_p"
$ws.Range("M1").Value = "network_param_m"
$ws.Range("N1").Value = "network_param_threshold"
$ws.Range("O1").Value = "vaccination_trust_percentage"
$ws.Range("P1").Value = "vaccination_ad_percentage"
$ws.Range("Q1").Value = "vaccination_ad_success_prob"
$ws.Range("R1").Value = "vaccination_action_prob"

# 3) New scenario row (row 2) values for the same columns.
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 0.2
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 0.6
$ws.Range("O2").Value = 0.5
$ws.Range("P2").Value = 0.1
$ws.Range("Q2").Value = 0.5
$ws.Range("R2").Value = 0.2

# 4) Match the header styling (the workbook's "label" font) used by all the
#    other header/parameter cells, by copying the format from the existing
#    H1:H2 pair (style already applied there) onto the new J:R cells.
$ws.Range("H1:H2").Copy()
$ws.Range("J1:R2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 5) Column widths for the new columns — match the widths used elsewhere in
#    the sheet for similarly-sized header text.
$ws.Range($ws.Columns.Item(8), $ws.Columns.Item(14)).ColumnWidth = 22.67
$ws.Range($ws.Columns.Item(15), $ws.Columns.Item(18)).ColumnWidth = 23.09

# 6) View tweaks captured in the saved file: zoom level and active selection.
$ws.Application.ActiveWindow.Zoom = 181
$null = $ws.Range("M2").Select()
